$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.103.04"
$ws.Range("E2").Value = "  +5.69%  "

$ws.Range("D3").Value = "1.922.08"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.85%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.67"
$ws.Range("E5").Value = "  +4.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5251"
$ws.Range("E7").Value = "  +3.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4062"
$ws.Range("E8").Value = "  +4.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08536"
$ws.Range("E9").Value = "  +2.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.10"
$ws.Range("E10").Value = "  +3.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.130"
$ws.Range("E11").Value = "  +2.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.58"
$ws.Range("E12").Value = "  +10.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.423"

$ws.Range("D14").Value = "1.919.25"
$ws.Range("E14").Value = "  +2.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.416"
$ws.Range("E15").Value = "  +2.07%  "

$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.74"
$ws.Range("E17").Value = "  +6.11%  "

$ws.Range("E18").Value = "  +1.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06709"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("E20").Value = "  +3.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.064"
$ws.Range("E22").Value = "  +2.76%  "

$ws.Range("D23").Value = "30.105.64"
$ws.Range("E23").Value = "  +5.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("E24").Value = "  +1.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.224"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("D26").Value = "2.140.61"
$ws.Range("E26").Value = "  +2.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.18"
$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.50"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.473"
$ws.Range("E29").Value = "  +3.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.73"
$ws.Range("E30").Value = "  +3.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.083"
$ws.Range("E31").Value = "  +4.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1060"
$ws.Range("E32").Value = "  +1.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.132"
$ws.Range("E33").Value = "  +6.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.649"
$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02526"
$ws.Range("E35").Value = "  +3.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06606"
$ws.Range("E36").Value = "  +1.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2233"
$ws.Range("E37").Value = "  +3.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.112"
$ws.Range("E38").Value = "  +3.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.239"
$ws.Range("E39").Value = "  +4.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.216"
$ws.Range("E40").Value = "  +3.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6572"
$ws.Range("E41").Value = "  +3.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.70"
$ws.Range("E42").Value = "  +5.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.246"
$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6205"
$ws.Range("E44").Value = "  +3.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.38"
$ws.Range("E45").Value = "  +2.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.788"
$ws.Range("E46").Value = "  +2.82%  "

$ws.Range("E47").Value = "  +4.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.249"
$ws.Range("E48").Value = "  +2.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.33"
$ws.Range("E49").Value = "  +2.90%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "80.07"
$ws.Range("E50").Value = "  +4.93%  "

$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.156"
$ws.Range("E51").Value = "  +1.05%  "
